function Set-TNRFont($rng) {
    $rng.Font.NameAscii = "Times New Roman"
    $rng.Font.NameOther = "Times New Roman"
    $rng.Font.NameBi = "Times New Roman"
    $rng.Font.Size = 12
    $rng.Font.SizeBi = 12
}

$d = $word.ActiveDocument

# --- Change 1: "Vamos ha utilizar el " -> "Vamos a utilizar el " ---
$d.Content.Find.Execute("Vamos ha utilizar el ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Vamos a utilizar el ", 2) | Out-Null

# --- Change 2: merge 4 runs ("Tratar las " / "imágenes" / " por " / "Patrones consiste...") into a single run ---
$p6 = $d.Paragraphs.Item(6)
$mergedText = "Tratar las imágenes por Patrones consiste procesar las imágenes y pasarlas la codificación que haría una calculadora de donde se marcan las líneas de un 8, en función del número que sea. Al ser letras la codificación sería:"
$delRng = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$delRng.Delete()
$p6b = $d.Paragraphs.Item(6)
$ins = $d.Range($p6b.Range.Start, $p6b.Range.Start)
$ins.InsertAfter($mergedText)
$p6c = $d.Paragraphs.Item(6)
Set-TNRFont($p6c.Range)

# --- Change 3: fill the empty paragraph with "Además, ..." text ---
$p11 = $d.Paragraphs.Item(11)
$p11.Range.InsertAfter("Además, se ha pensado que sería buena idea recortar los bordes de las imágenes antes de tratarlas. Sin embargo, también se tratarán imágenes sin recortar.")
$p11b = $d.Paragraphs.Item(11)
Set-TNRFont($p11b.Range)

# --- Change 4: add trailing space to " y Árboles de decisión." ---
$d.Content.Find.Execute(" y Árboles de decisión.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " y Árboles de decisión. ", 2) | Out-Null

# --- Change 5: remove the old _GoBack bookmark + trailing " " run from paragraph 12 ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$p12 = $d.Paragraphs.Item(12)
$trailSpaceRng = $d.Range($p12.Range.End - 2, $p12.Range.End - 1)
$trailSpaceRng.Delete()

# --- Change 5b: insert the new tail paragraphs via raw OOXML ---
$p12b = $d.Paragraphs.Item(12)
$insertPt = $d.Range($p12b.Range.End - 1, $p12b.Range.End - 1)

$tailXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Grupo 03:</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr>
    <w:contextualSpacing/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Daniel Fernández</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:contextualSpacing/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Román García</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:contextualSpacing/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Patricia Losana</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:contextualSpacing/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="26"/>
      <w:szCs w:val="26"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Celia San Gregorio</w:t>
  </w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPt.InsertXML($tailXml)
